$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")

# Find the row containing "The Legacy Virus (Custom Achievement)" in column A and delete it entirely.
$found = $ws.Range("A1:A100").Find("The Legacy Virus (Custom Achievement)")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}
